# 03 classification validation & restructure future labs
# Regenerate the "Payment Cards" sample data (rows 2-13) with a new
# fictitious data set for the updated classification-validation lab,
# and tighten the Email column (F) to its new best-fit width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer ID, First Name, Last Name, Credit Card Number, Credit Card Type, Email
$rows = @(
    @('C001', 'Lauren',   'Smith',     '6011-3336-6513-2206', 'Discover',   'lauren.smith@contoso.com'),
    @('C002', 'Robert',   'Rodriguez', '3747-147727-40130',   'Amex',       'robert.rodriguez@contoso.com'),
    @('C003', 'Lauren',   'Smith',     '3705-465003-28226',   'Amex',       'lauren.smith@gmail.com'),
    @('C004', 'Amanda',   'Wilson',    '5124-4082-4524-2787', 'Mastercard', 'amanda.wilson@contoso.com'),
    @('C005', 'Jessica',  'Moore',     '4844-4614-3142-2578', 'Visa',       'jessica.moore@contoso.com'),
    @('C006', 'Lauren',   'Davis',     '5203-4474-0703-1417', 'Mastercard', 'lauren.davis@hotmail.com'),
    @('C007', 'Emily',    'Davis',     '3704-733025-40805',   'Amex',       'emily.davis@contoso.com'),
    @('C008', 'Jennifer', 'Gonzalez',  '4565-4310-5433-1647', 'Visa',       'jennifer.gonzalez@contoso.com'),
    @('C009', 'Michael',  'Taylor',    '6011-7227-1566-3009', 'Discover',   'michael.taylor@icloud.com'),
    @('C010', 'Michelle', 'Taylor',    '5237-0072-8150-0588', 'Mastercard', 'michelle.taylor@contoso.com'),
    @('C011', 'Sarah',    'Wilson',    '4271-0271-1283-3614', 'Visa',       'sarah.wilson@contoso.com'),
    @('C012', 'James',    'Garcia',    '6011-8345-0214-0321', 'Discover',   'james.garcia@hotmail.com')
)

$columns = @('A', 'B', 'C', 'D', 'E', 'F')
$startRow = 2

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $addr = "$($columns[$j])$rowNum"
        $ws.Range($addr).Value = $rowData[$j]
    }
}

# New Email values are shorter than before; shrink column F's best-fit width to match (31.28... -> 30).
$ws.Columns.Item(6).ColumnWidth = 29.15

Write-Output "Updated $($rows.Count) data rows and resized column F."
